$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.212.23"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "3.090.82"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.89"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.082.72"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.61"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.158"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.453"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "34.95"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000216"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "3.585.06"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "63.318.31"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "3.093.18"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "502.48"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.65"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.46"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.703"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.23"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.50"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.20"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.75"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.14"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.97"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.17"
$ws.Range("D31").ClearFormats()
$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.11"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.50"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.68"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +13.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "526.12"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -8.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.86"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.15"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0410"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.121"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0787"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").Value = "3.042.66"
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.05"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.59"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.95%  "
$ws.Range("E44").Value = "  +4.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.05"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.34"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.81%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.65"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.77%  "
$ws.Range("E50").Value = "  -3.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.32"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +67.47%  "
